$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.808608104485756
$ws.Range("C2").Value = 0.1215340199039758
$ws.Range("D2").Value = 0.4710237188118214
$ws.Range("E2").Value = 0.1400167628055051
$ws.Range("G2").Value = 0.002507324897131591
$ws.Range("I2").Value = 1.169100895398643
$ws.Range("J2").Value = 0.0470518017925956
$ws.Range("L2").Value = 0.5266700258374613
$ws.Range("N2").Value = 1.548036878411128
$ws.Range("O2").Value = 5.431941790156287

$ws.Range("B3").Value = 1.711197344110815
$ws.Range("C3").Value = 0.1062177000250415
$ws.Range("D3").Value = 0.4699262764096659
$ws.Range("E3").Value = 0.1406612279410435
$ws.Range("G3").Value = 0.002510735781739136
$ws.Range("I3").Value = 1.177282801318235
$ws.Range("J3").Value = 0.04703706208949043
$ws.Range("L3").Value = 0.5185733965032995
$ws.Range("N3").Value = 1.563432616638543
$ws.Range("O3").Value = 5.44181348034482

$ws.Range("B4").Value = 1.651937334414754
$ws.Range("C4").Value = 0.09677480672948491
$ws.Range("D4").Value = 0.4694543860747729
$ws.Range("E4").Value = 0.1411037367504324
$ws.Range("G4").Value = 0.002512943342198723
$ws.Range("I4").Value = 1.182909629404168
$ws.Range("J4").Value = 0.04703029930473868
$ws.Range("L4").Value = 0.5138232423198588
$ws.Range("N4").Value = 1.573457637170023
$ws.Range("O4").Value = 5.450854641259866

$ws.Range("B5").Value = 1.627928331145824
$ws.Range("C5").Value = 0.09291716556785445
$ws.Range("D5").Value = 0.4693129397010409
$ws.Range("E5").Value = 0.1412958493672036
$ws.Range("G5").Value = 0.002513871509479405
$ws.Range("I5").Value = 1.185354194434854
$ws.Range("J5").Value = 0.04702812166982717
$ws.Range("L5").Value = 0.5119432929360954
$ws.Range("N5").Value = 1.577686807104921
$ws.Range("O5").Value = 5.455287987682027

$ws.Range("B6").Value = 1.623950149726454
$ws.Range("C6").Value = 0.09227603191246203
$ws.Range("D6").Value = 0.4692925261604444
$ws.Range("E6").Value = 0.1413284619172916
$ws.Range("G6").Value = 0.00251402735912021
$ws.Range("I6").Value = 1.185769266085995
$ws.Range("J6").Value = 0.04702779509026289
$ws.Range("L6").Value = 0.5116345020693984
$ws.Range("N6").Value = 1.578397750310209
$ws.Range("O6").Value = 5.456069370012614

$ws.Range("B7").Value = 1.651612972057137
$ws.Range("C7").Value = 0.09672281987724318
$ws.Range("D7").Value = 0.4694522724813623
$ws.Range("E7").Value = 0.1411062799014076
$ws.Range("G7").Value = 0.002512955744097271
$ws.Range("I7").Value = 1.182941983987256
$ws.Range("J7").Value = 0.04703026759076501
$ws.Range("L7").Value = 0.5137976626225935
$ws.Range("N7").Value = 1.573514090532914
$ws.Range("O7").Value = 5.450911398685349

$ws.Range("B8").Value = 1.774907453174421
$ws.Range("C8").Value = 0.1162610360283338
$ws.Range("D8").Value = 0.4706034491445905
$ws.Range("E8").Value = 0.1402292703709183
$ws.Range("G8").Value = 0.002508477514925445
$ws.Range("I8").Value = 1.171796799938917
$ws.Range("J8").Value = 0.04704624659344159
$ws.Range("L8").Value = 0.5238324835700183
$ws.Range("N8").Value = 1.553226654520412
$ws.Range("O8").Value = 5.434726863279991

$ws.Range("B9").Value = 2.021004588198991
$ws.Range("C9").Value = 0.1542654390095777
$ws.Range("D9").Value = 0.4744606917072076
$ws.Range("E9").Value = 0.1388801068182719
$ws.Range("G9").Value = 0.002500590432057225
$ws.Range("I9").Value = 1.154730065022498
$ws.Range("J9").Value = 0.04709560435760807
$ws.Range("L9").Value = 0.545260796720811
$ws.Range("N9").Value = 1.517978314925323
$ws.Range("O9").Value = 5.426653879013571

$ws.Range("B10").Value = 2.204394298454588
$ws.Range("C10").Value = 0.1819961124008103
$ws.Range("D10").Value = 0.4782674024891662
$ws.Range("E10").Value = 0.1381139129658866
$ws.Range("G10").Value = 0.002495335653980074
$ws.Range("I10").Value = 1.145116230419383
$ws.Range("J10").Value = 0.0471426949923579
$ws.Range("L10").Value = 0.5620664579082302
$ws.Range("N10").Value = 1.494841628519957
$ws.Range("O10").Value = 5.435186879360487

$ws.Range("B11").Value = 2.288374492314006
$ws.Range("C11").Value = 0.1945697096198273
$ws.Range("D11").Value = 0.4802099514723466
$ws.Range("E11").Value = 0.1378140269634578
$ws.Range("G11").Value = 0.002493061158689865
$ws.Range("I11").Value = 1.141379096156086
$ws.Range("J11").Value = 0.04716643586013802
$ws.Range("L11").Value = 0.5699416348258097
$ws.Range("N11").Value = 1.484914637665305
$ws.Range("O11").Value = 5.442218311051732

$ws.Range("B12").Value = 2.320254251633173
$ws.Range("C12").Value = 0.1993249888973594
$ws.Range("D12").Value = 0.4809758100278287
$ws.Range("E12").Value = 0.1377074493022086
$ws.Range("G12").Value = 0.002492216447189879
$ws.Range("I12").Value = 1.140055534883778
$ws.Range("J12").Value = 0.04717575648440686
$ws.Range("L12").Value = 0.5729567457151887
$ws.Range("N12").Value = 1.481241499302989
$ws.Range("O12").Value = 5.44533443177437

$ws.Range("B13").Value = 2.31338491832588
$ws.Range("C13").Value = 0.1983011265865571
$ws.Range("D13").Value = 0.4808095239869346
$ws.Range("E13").Value = 0.1377300923826219
$ws.Range("G13").Value = 0.002492397634357886
$ws.Range("I13").Value = 1.140336511186028
$ws.Range("J13").Value = 0.04717373447116557
$ws.Range("L13").Value = 0.5723059244524649
$ws.Range("N13").Value = 1.482028750459556
$ws.Range("O13").Value = 5.444643142508482

$ws.Range("B14").Value = 2.290995699362668
$ws.Range("C14").Value = 0.194961051748777
$ws.Range("D14").Value = 0.4802723530755912
$ws.Range("E14").Value = 0.1378051188948017
$ws.Range("G14").Value = 0.002492991331859079
$ws.Range("I14").Value = 1.141268369281491
$ws.Range("J14").Value = 0.04716719606739517
$ws.Range("L14").Value = 0.5701890302552073
$ws.Range("N14").Value = 1.484610722791217
$ws.Range("O14").Value = 5.442465584833343

$ws.Range("B15").Value = 2.277291804532183
$ws.Range("C15").Value = 0.1929143649151115
$ws.Range("D15").Value = 0.4799472589108262
$ws.Range("E15").Value = 0.1378519837198784
$ws.Range("G15").Value = 0.002493357146947846
$ws.Range("I15").Value = 1.141851093188876
$ws.Range("J15").Value = 0.04716323405467548
$ws.Range("L15").Value = 0.5688966583662705
$ws.Range("N15").Value = 1.48620345507787
$ws.Range("O15").Value = 5.441190839286008

$ws.Range("B16").Value = 2.198917035031343
$ws.Range("C16").Value = 0.1811735566337802
$ws.Range("D16").Value = 0.4781446898853261
$ws.Range("E16").Value = 0.1381344891306551
$ws.Range("G16").Value = 0.002495486623512016
$ws.Range("I16").Value = 1.145373273983672
$ws.Range("J16").Value = 0.04714118987952354
$ws.Range("L16").Value = 0.5615564159121362
$ws.Range("N16").Value = 1.495502414425452
$ws.Range("O16").Value = 5.434790785645504

$ws.Range("B17").Value = 2.150977758977376
$ws.Range("C17").Value = 0.1739602918065941
$ws.Range("D17").Value = 0.4770928323730175
$ws.Range("E17").Value = 0.1383202495360116
$ws.Range("G17").Value = 0.002496822623739309
$ws.Range("I17").Value = 1.147697060195178
$ws.Range("J17").Value = 0.04712825832105949
$ws.Range("L17").Value = 0.5571122671440918
$ws.Range("N17").Value = 1.501360190089713
$ws.Range("O17").Value = 5.431671665655358

$ws.Range("B18").Value = 2.123456715635939
$ws.Range("C18").Value = 0.1698075407483941
$ws.Range("D18").Value = 0.4765076824037919
$ws.Range("E18").Value = 0.1384316749234813
$ws.Range("G18").Value = 0.002497601972281973
$ws.Range("I18").Value = 1.149093519990402
$ws.Range("J18").Value = 0.0471210388942076
$ws.Range("L18").Value = 0.5545777864876555
$ws.Range("N18").Value = 1.504785713996071
$ws.Range("O18").Value = 5.430174064708524

$ws.Range("B19").Value = 2.114147604065295
$ws.Range("C19").Value = 0.1684008298939261
$ws.Range("D19").Value = 0.476312971990879
$ws.Range("E19").Value = 0.1384701888656501
$ws.Range("G19").Value = 0.002497867723748711
$ws.Range("I19").Value = 1.149576618035418
$ws.Range("J19").Value = 0.04711863213370293
$ws.Range("L19").Value = 0.5537233823657743
$ws.Range("N19").Value = 1.505955204681591
$ws.Range("O19").Value = 5.429717899201194

$ws.Range("B20").Value = 2.156075569034442
$ws.Range("C20").Value = 0.1747285579135678
$ws.Range("D20").Value = 0.4772027503954348
$ws.Range("E20").Value = 0.1383000010217579
$ws.Range("G20").Value = 0.002496679274958116
$ws.Range("I20").Value = 1.147443491029208
$ws.Range("J20").Value = 0.04712961231901502
$ws.Range("L20").Value = 0.5575831120533081
$ws.Range("N20").Value = 1.500730794425536
$ws.Range("O20").Value = 5.431973017522381

$ws.Range("B21").Value = 2.297569843594545
$ws.Range("C21").Value = 0.1959422785270988
$ws.Range("D21").Value = 0.4804293125209966
$ws.Range("E21").Value = 0.1377828923971656
$ws.Range("G21").Value = 0.002492816498823184
$ws.Range("I21").Value = 1.140992172739104
$ws.Range("J21").Value = 0.0471691076101397
$ws.Range("L21").Value = 0.5708099202687436
$ws.Range("N21").Value = 1.483850000909051
$ws.Range("O21").Value = 5.443092874903073

$ws.Range("B22").Value = 2.390500098755581
$ws.Range("C22").Value = 0.2097712308132316
$ws.Range("D22").Value = 0.4827143915480434
$ws.Range("E22").Value = 0.137485625653964
$ws.Range("G22").Value = 0.002490388616194934
$ws.Range("I22").Value = 1.137309894520762
$ws.Range("J22").Value = 0.04719684493131737
$ws.Range("L22").Value = 0.5796464063604674
$ws.Range("N22").Value = 1.47331869275181
$ws.Range("O22").Value = 5.453003771253748

$ws.Range("B23").Value = 2.340860304678586
$ws.Range("C23").Value = 0.2023937479658571
$ws.Range("D23").Value = 0.4814786881792656
$ws.Range("E23").Value = 0.1376405638007547
$ws.Range("G23").Value = 0.002491675604198463
$ws.Range("I23").Value = 1.139226292291681
$ws.Range("J23").Value = 0.04718186586557316
$ws.Range("L23").Value = 0.5749126884894196
$ws.Range("N23").Value = 1.478893581070317
$ws.Range("O23").Value = 5.44747208484813

$ws.Range("B24").Value = 2.153770724744618
$ws.Range("C24").Value = 0.1743812427053513
$ws.Range("D24").Value = 0.4771529954736167
$ws.Range("E24").Value = 0.1383091409544956
$ws.Range("G24").Value = 0.002496744047781285
$ws.Range("I24").Value = 1.14755794124833
$ws.Range("J24").Value = 0.04712899950645699
$ws.Range("L24").Value = 0.5573701791436321
$ws.Range("N24").Value = 1.501015164122549
$ws.Range("O24").Value = 5.431835855515146

$ws.Range("B25").Value = 1.953971780879442
$ws.Range("C25").Value = 0.1440176548362047
$ws.Range("D25").Value = 0.4732461949940756
$ws.Range("E25").Value = 0.1392055064252062
$ws.Range("G25").Value = 0.002502628890778834
$ws.Range("I25").Value = 1.158833802796643
$ws.Range("J25").Value = 0.04708033887839136
$ws.Range("L25").Value = 0.5392769934057497
$ws.Range("N25").Value = 1.527028884861124
$ws.Range("O25").Value = 5.426300495964568
